$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 483, shifting existing rows 483-530 down to 484-531
$ws.Rows("483:483").Insert()

# Fill in the new row 483 with the inserted record's data
$ws.Cells.Item(483, 1).Value = 4
$ws.Cells.Item(483, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(483, 3).Value = "Los Lagos"
$ws.Cells.Item(483, 4).Value = 45212
$ws.Cells.Item(483, 5).Value = 10
$ws.Cells.Item(483, 6).Value = 100112037
$ws.Cells.Item(483, 7).Value = "Cebollín"
$ws.Cells.Item(483, 8).Value = "Sin especificar"
$ws.Cells.Item(483, 9).Value = "Primera"
$ws.Cells.Item(483, 10).Value = 180
$ws.Cells.Item(483, 11).Value = 6500
$ws.Cells.Item(483, 12).Value = 6500
$ws.Cells.Item(483, 13).Value = 6500
$ws.Cells.Item(483, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(483, 15).Value = "Región Metropolitana"
$ws.Cells.Item(483, 16).Value = 181
$ws.Cells.Item(483, 17).Value = 36
$ws.Cells.Item(483, 18).Value = "Hortaliza"
